$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Parameter..."/"Values..." columns
# (AR/AS) so they shift right to AS/AT, leaving three fresh blank columns
# (AO, AP, AQ) right after the existing "goal_name_id" column (AN).
$ws.Columns("AR").Insert()

# Copy the header formatting of AN1 ("goal_name_id") onto the three new
# header cells, then fill in the new header labels.
$ws.Range("AN1").Copy()
$ws.Range("AO1:AQ1").PasteSpecial(-4122)
$ws.Range("AO1").Value = "goal_version"
$ws.Range("AP1").Value = "rule_name_id"
$ws.Range("AQ1").Value = "rule_version"

# Copy the data-row formatting of AN2 onto the new AO2 cell (stays blank,
# same as the other placeholder cells in that row).
$ws.Range("AN2").Copy()
$ws.Range("AO2").PasteSpecial(-4122)

# Match the "Parameter..." column's width on the newly inserted AP column
# too (the raw stored width is ColumnWidth + 0.8333333333333).
$ws.Range("AP1").ColumnWidth = 14.1666666666667

# Update the view: scroll toward the new columns and select AN2, matching
# the post-edit selection state.
[void]$ws.Range("AN2").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 34
$win.ScrollRow = 1
